# Auto-generated edit script: updates market price / profit figures
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM) per refreshed
# Universalis market data snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 164839
$ws.Range("I6").Value = 246808.5
$ws.Range("K6").Value = 740425.5
$ws.Range("M6").Value = -740313.5

# Row 62
$ws.Range("H62").Value = 2686.5557
$ws.Range("J62").Value = 2882.7144
$ws.Range("L62").Value = 2882.7144
$ws.Range("N62").Value = -4130.7144

# Row 65
$ws.Range("H65").Value = 2686.5557
$ws.Range("J65").Value = 2882.7144
$ws.Range("L65").Value = 14413.572
$ws.Range("N65").Value = -20653.572

# Row 98
$ws.Range("H98").Value = 425
$ws.Range("I98").Value = 420.58823
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 420.58823
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 1077.41177
$ws.Range("N98").Value = -3496

# Row 121
$ws.Range("H121").Value = 1165.8889
$ws.Range("J121").Value = 1082.1666
$ws.Range("L121").Value = 3246.4998
$ws.Range("N121").Value = -6740.4998

# Row 122
$ws.Range("H122").Value = 425
$ws.Range("I122").Value = 420.58823
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 1261.76469
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = 1188.23531
$ws.Range("N122").Value = -6400

# Row 129
$ws.Range("H129").Value = 950.22
$ws.Range("I129").Value = 544.4706
$ws.Range("J129").Value = 1033.3253
$ws.Range("K129").Value = 1633.4118
$ws.Range("L129").Value = 3099.9759
$ws.Range("M129").Value = 3366.5882
$ws.Range("N129").Value = -13099.9759

# Row 141
$ws.Range("H141").Value = 2433.3333
$ws.Range("I141").Value = 2173.3333
$ws.Range("J141").Value = 3733.3333
$ws.Range("K141").Value = 6519.999899999999
$ws.Range("L141").Value = 11199.9999
$ws.Range("M141").Value = -1339.999899999999
$ws.Range("N141").Value = -21559.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1507.0526
$ws.Range("I61").Value = 1377.125
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 1377.125
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -1165.125
$ws.Range("N61").Value = -2624

# Row 74
$ws.Range("H74").Value = 1092.7778
$ws.Range("I74").Value = 1092.7778
$ws.Range("K74").Value = 1092.7778
$ws.Range("M74").Value = -218.7778000000001

# Row 77
$ws.Range("H77").Value = 1092.7778
$ws.Range("I77").Value = 1092.7778
$ws.Range("K77").Value = 5463.889
$ws.Range("M77").Value = -1095.889

# Row 136
$ws.Range("H136").Value = 1507.0526
$ws.Range("I136").Value = 1377.125
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 4131.375
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -1581.375
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 71622.25
$ws.Range("I86").Value = 95017.164
$ws.Range("J86").Value = 1437.5
$ws.Range("K86").Value = 95017.164
$ws.Range("L86").Value = 1437.5
$ws.Range("M86").Value = -93894.164
$ws.Range("N86").Value = -3683.5

# Row 89
$ws.Range("H89").Value = 71622.25
$ws.Range("I89").Value = 95017.164
$ws.Range("J89").Value = 1437.5
$ws.Range("K89").Value = 475085.82
$ws.Range("L89").Value = 7187.5
$ws.Range("M89").Value = -469469.82
$ws.Range("N89").Value = -18419.5

$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 66734.57000000001
$ws.Range("J9").Value = 66734.57000000001
$ws.Range("L9").Value = 66734.57000000001
$ws.Range("N9").Value = -67070.57000000001

# Row 31
$ws.Range("H31").Value = 22154.2
$ws.Range("I31").Value = 1431.091
$ws.Range("J31").Value = 40636.973
$ws.Range("K31").Value = 1431.091
$ws.Range("L31").Value = 40636.973
$ws.Range("M31").Value = -1136.091
$ws.Range("N31").Value = -41226.973

# Row 34
$ws.Range("H34").Value = 22154.2
$ws.Range("I34").Value = 1431.091
$ws.Range("J34").Value = 40636.973
$ws.Range("K34").Value = 1431.091
$ws.Range("L34").Value = 40636.973
$ws.Range("M34").Value = -1229.091
$ws.Range("N34").Value = -41040.973

# Row 99
$ws.Range("H99").Value = 14439.223
$ws.Range("I99").Value = 5095
$ws.Range("J99").Value = 21914.6
$ws.Range("K99").Value = 5095
$ws.Range("L99").Value = 21914.6
$ws.Range("M99").Value = -3597
$ws.Range("N99").Value = -24910.6

# Row 122
$ws.Range("H122").Value = 1953.625
$ws.Range("I122").Value = 850
$ws.Range("J122").Value = 2321.5
$ws.Range("K122").Value = 2550
$ws.Range("L122").Value = 6964.5
$ws.Range("M122").Value = -100
$ws.Range("N122").Value = -11864.5

# Row 126
$ws.Range("H126").Value = 14439.223
$ws.Range("I126").Value = 5095
$ws.Range("J126").Value = 21914.6
$ws.Range("K126").Value = 15285
$ws.Range("L126").Value = 65743.79999999999
$ws.Range("M126").Value = -12815
$ws.Range("N126").Value = -70683.79999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 1930
$ws.Range("I3").Value = 1930
$ws.Range("K3").Value = 5790
$ws.Range("M3").Value = -5678

# Row 5
$ws.Range("H5").Value = 1358.5834
$ws.Range("J5").Value = 1474.3478
$ws.Range("L5").Value = 4423.0434
$ws.Range("N5").Value = -4647.0434

# Row 97
$ws.Range("H97").Value = 1333.7778
$ws.Range("I97").Value = 1340
$ws.Range("K97").Value = 4020
$ws.Range("M97").Value = -3524

# Row 131
$ws.Range("H131").Value = 791.09
$ws.Range("J131").Value = 818.8421
$ws.Range("L131").Value = 2456.5263
$ws.Range("N131").Value = -12536.5263

# Row 135
$ws.Range("H135").Value = 1358.5834
$ws.Range("J135").Value = 1474.3478
$ws.Range("L135").Value = 13269.1302
$ws.Range("N135").Value = -18339.1302

# Row 138
$ws.Range("H138").Value = 9926.77
$ws.Range("I138").Value = 13873.5
$ws.Range("J138").Value = 3612
$ws.Range("K138").Value = 41620.5
$ws.Range("L138").Value = 10836
$ws.Range("M138").Value = -36480.5
$ws.Range("N138").Value = -21116

# Row 139
$ws.Range("H139").Value = 1698.0454
$ws.Range("I139").Value = 787.4375
$ws.Range("J139").Value = 4126.3335
$ws.Range("K139").Value = 2362.3125
$ws.Range("L139").Value = 12379.0005
$ws.Range("M139").Value = 2777.6875
$ws.Range("N139").Value = -22659.0005

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 43447.81
$ws.Range("I70").Value = 71390
$ws.Range("K70").Value = 71390
$ws.Range("M70").Value = -71120

# Row 73
$ws.Range("H73").Value = 43447.81
$ws.Range("I73").Value = 71390
$ws.Range("K73").Value = 71390
$ws.Range("M73").Value = -70454
